$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

# Create a fresh paragraph right after it, styled as a bullet list item.
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Style = "ListBullet"

# Replace the whole (still-empty) paragraph -- including its paragraph
# mark -- with the three professor lines, each its own run, separated by
# manual line breaks (matches how the rest of the document represents
# multi-line list items).
$rng = $newPara.Range

$runsXml = (
    '<w:r><w:t>7459752 - Maria Ismenia Sodero Toledo Faria</w:t><w:br/></w:r>' +
    '<w:r><w:t>2166002 - Sandra Giacomin Schneider</w:t><w:br/></w:r>' +
    '<w:r><w:t>1922320 - Sebastiao Ribeiro</w:t></w:r>'
)

$paraXml = '<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' + $runsXml + '</w:p>'

$packageXml = (
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $paraXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
)

$rng.InsertXML($packageXml)
